# Apply the facebook_groups_report.xlsx edits to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Update the group-keywords text in C23 (shared string content change).
$ws.Range("C23").Value = "BITCOIN, BNB, METAVERSE, BLOCKCHAIN, TRON, OPENSEA, NFT"

# 2. Update the "Number of Members" (E) and "Number of Posts from Outside" (F)
#    columns for rows 2-29 with the refreshed report numbers.
$updates = @{
    2  = @{ E = 17035; F = 23 }
    3  = @{ E = 1527;  F = 6  }
    4  = @{ E = 16211; F = 30 }
    5  = @{ E = 16182; F = 20 }
    6  = @{ E = 16128; F = 24 }
    7  = @{ E = 16094; F = 17 }
    8  = @{ E = 15888; F = 8  }
    9  = @{ E = 15843; F = 13 }
    10 = @{ E = 15762; F = 27 }
    11 = @{ E = 15686; F = 17 }
    12 = @{ E = 1457;  F = 4  }
    13 = @{ E = 15129; F = 8  }
    14 = @{ E = 14003; F = 11 }
    15 = @{ E = 13930; F = 6  }
    16 = @{ E = 13871; F = 7  }
    17 = @{ E = 12534; F = 5  }
    18 = @{ E = 13414; F = 9  }
    19 = @{ E = 1307;  F = 3  }
    20 = @{ E = 13303; F = 6  }
    21 = @{ E = 13211; F = 1  }
    22 = @{ E = 13165; F = 38 }
    23 = @{ E = 1298;  F = 1  }
    24 = @{ E = 16570; F = 11 }
    25 = @{ E = 16411; F = 5  }
    26 = @{ E = 1511;  F = 2  }
    27 = @{ E = 16310; F = 19 }
    28 = @{ E = 16124; F = 5  }
    29 = @{ E = 16200; F = 4  }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 5).Value = $vals.E
    $ws.Cells.Item($row, 6).Value = $vals.F
}

# 3. Rows 30-33: refresh E, and collapse the F "rollup" SUM formulas down to
#    a plain literal 0 (the subtree they summed no longer exists).
$rollupUpdates = @{
    30 = 398
    31 = 3534
    32 = 3514
    33 = 3504
}
foreach ($row in $rollupUpdates.Keys) {
    $ws.Cells.Item($row, 5).Value = $rollupUpdates[$row]
    $ws.Cells.Item($row, 6).Value = 0
}

# Row 34: only E changes, F34 already a plain 0.
$ws.Cells.Item(34, 5).Value = 3492

# 4. Move the active cell/selection from H6 to B2.
$ws.Range("B2").Select() | Out-Null
